$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width adjustments (A, B, C) -----------------------------------
# Note: Excel's ColumnWidth setter quantizes internally, so the closest
# achievable width is used.
$ws.Columns.Item(1).ColumnWidth = 10.875
$ws.Columns.Item(2).ColumnWidth = 10.875
$ws.Columns.Item(3).ColumnWidth = 8.875

# --- EC consumption profile values (rows 6, 12, 18, 24, 30, 36, 42, 48, 54, 60) ---
$ws.Range("A6").Value = 6659.7100000000019
$ws.Range("B6").Value = 3682.746000000001
$ws.Range("C6").Value = 4209.5439999999999

$ws.Range("A12").Value = 6190.8819999999969
$ws.Range("B12").Value = 3433.9150000000004
$ws.Range("C12").Value = 4066.0949999999993

$ws.Range("A18").Value = 6418.5089999999991
$ws.Range("B18").Value = 3538.2750000000015
$ws.Range("C18").Value = 4137.1259999999975

$ws.Range("A24").Value = 6446.7100000000028
$ws.Range("B24").Value = 2889.52
$ws.Range("C24").Value = 3212.8600000000006

$ws.Range("A30").Value = 5620.2959999999975
$ws.Range("B30").Value = 3200.3049999999994
$ws.Range("C30").Value = 3448.8949999999995

$ws.Range("A36").Value = 5041.8729999999978
$ws.Range("B36").Value = 3122.7190000000005
$ws.Range("C36").Value = 3308.4050000000007

$ws.Range("A42").Value = 5868.4619999999932
$ws.Range("B42").Value = 3155.5130000000022
$ws.Range("C42").Value = 3470.2359999999999

$ws.Range("A48").Value = 5499.7400000000034
$ws.Range("B48").Value = 3429.6760000000004
$ws.Range("C48").Value = 3680.5840000000017

$ws.Range("A54").Value = 4869.9800000000023
$ws.Range("B54").Value = 2704.0450000000005
$ws.Range("C54").Value = 2941.7550000000001

$ws.Range("A60").Value = 4328.6750000000002
$ws.Range("B60").Value = 2233.5149999999994
$ws.Range("C60").Value = 2605.3999999999987
